# edit.ps1 — apply the template_recommandations_Mobile_Phone.docx changes:
#  1) Adjust several paragraphs' "space before" values.
#  2) Split the PIN line into 3 runs so "ne pas changer" is bold/red, and
#     reword it from "Num. PIN. (ne pas changer) : {{PIN}} " to
#     "PIN : {{PIN}} (ne pas changer) ".
#  3) Same treatment for the PUK line.
#  4) Rename the {{Used_by}} placeholder to {{Utilise_par}}.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Spacing-before tweaks (values are in points; OOXML w:before is in
#    twentieths of a point, so e.g. 567/20 = 28.35).
# ---------------------------------------------------------------------

# Title paragraph: before=796 -> before=567
$d.Paragraphs.Item(3).SpaceBefore = 567 / 20

# IMEI paragraph: before=588 -> before=283
$d.Paragraphs.Item(7).SpaceBefore = 283 / 20

# PIN paragraph: before=1124 -> before=283
$d.Paragraphs.Item(8).SpaceBefore = 283 / 20

# Consignes paragraph: before=416 -> before=850
$d.Paragraphs.Item(11).SpaceBefore = 850 / 20

# "Par la présente" paragraph: before=694 -> before=1417
$d.Paragraphs.Item(19).SpaceBefore = 1417 / 20

# ---------------------------------------------------------------------
# 2) PIN line: "Num. PIN. (ne pas changer) : {{PIN}} "
#           -> "PIN : {{PIN}} (" + "ne pas changer" (bold, red) + ") "
# ---------------------------------------------------------------------

$pPin = $d.Paragraphs.Item(8)
$pPin.Range.Text = "PIN : {{PIN}} (ne pas changer) "
$pinStart = $pPin.Range.Start
$pinEnd = $pPin.Range.End

$pinHi = $d.Range($pinStart, $pinEnd)
$pinHi.Find.Execute("ne pas changer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pinHi.Font.Bold = $true
$pinHi.Font.Color = 0x1E21C9

# ---------------------------------------------------------------------
# 3) PUK line: "Num. PUK. (ne pas changer) : {{PUK}} "
#           -> "PUK : {{PUK}} (" + "ne pas changer" (bold, red) + ") "
# ---------------------------------------------------------------------

$pPuk = $d.Paragraphs.Item(9)
$pPuk.Range.Text = "PUK : {{PUK}} (ne pas changer) "
$pukStart = $pPuk.Range.Start
$pukEnd = $pPuk.Range.End

$pukHi = $d.Range($pukStart, $pukEnd)
$pukHi.Find.Execute("ne pas changer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pukHi.Font.Bold = $true
$pukHi.Font.Color = 0x1E21C9

# ---------------------------------------------------------------------
# 4) Placeholder rename: {{Used_by}} -> {{Utilise_par}}
# ---------------------------------------------------------------------

$d.Content.Find.Execute("Used_by", $true, $false, $false, $false, $false, $true, 1, $false, "Utilise_par", 2)

Write-Output "done"
